$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.660.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.742.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.04%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'333.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.17%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.3762"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.32%  "
$ws.Range("D8").Value = "'48.49"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.35%  "
$ws.Range("D9").Value = "'0.3381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.86%  "
$ws.Range("D10").Value = "'1.186"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.84%  "
$ws.Range("D11").Value = "'0.07465"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.70%  "
$ws.Range("D12").Value = "'1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'6.435"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.42%  "
$ws.Range("D14").Value = "'20.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.58%  "
$ws.Range("D15").Value = "'7.127"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.25%  "
$ws.Range("D16").Value = "'1.739.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.15%  "
$ws.Range("D17").Value = "'0.00001085"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.84%  "
$ws.Range("D18").Value = "'0.06675"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "'83.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.88%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'16.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.27%  "
$ws.Range("D22").Value = "'6.193"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.03%  "
$ws.Range("D23").Value = "'12.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("D24").Value = "'26.603.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.06%  "
$ws.Range("D25").Value = "'2.452"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "'2.442"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("D27").Value = "'1.404"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +17.02%  "
$ws.Range("D28").Value = "'153.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.32%  "
$ws.Range("D29").Value = "'19.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.21%  "
$ws.Range("D30").Value = "'1.934.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.34%  "
$ws.Range("D31").Value = "'131.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.12%  "
$ws.Range("D32").Value = "'4.150"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").Value = "'6.096"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.79%  "
$ws.Range("D34").Value = "'0.08635"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("D35").Value = "'1.707"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.32%  "
$ws.Range("D36").Value = "'12.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.91%  "
$ws.Range("D37").Value = "'5.420"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.01%  "
$ws.Range("D38").Value = "'0.02356"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.58%  "
$ws.Range("D39").Value = "'0.06296"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.01%  "
$ws.Range("D40").Value = "'0.2179"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.64%  "
$ws.Range("D41").Value = "'8.594"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.65%  "
$ws.Range("E42").Value = "  -4.33%  "
$ws.Range("D43").Value = "'0.6229"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.14%  "
$ws.Range("D44").Value = "'14.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.57%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "'3.910"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.28%  "
$ws.Range("D47").Value = "'0.6046"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.79%  "
$ws.Range("D48").Value = "'128.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").Value = "'2.057"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.95%  "
$ws.Range("D50").Value = "'0.07241"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("D51").Value = "'77.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.25%  "
